# Fruta / hortaliza, semanal
# Insert two new weekly records (rows 49-50) for the "Damasco" sheet, pushing the
# previously-existing rows 49-51 down to become rows 51-53 (their content is
# unchanged, Excel's row-insert naturally preserves it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new blank rows at position 49 (shifts old rows 49,50,51 -> 51,52,53)
$ws.Rows("49:50").Insert()

# --- New row 49: Damasco, Modesto, Especial, Región de O'Higgins, 18kg box ---
$ws.Cells.Item(49, 1).Value = 8
$ws.Cells.Item(49, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(49, 3).Value = "Coquimbo"
$ws.Cells.Item(49, 4).Value = 44568
$ws.Cells.Item(49, 5).Value = 4
$ws.Cells.Item(49, 6).Value = "Fruta"
$ws.Cells.Item(49, 7).Value = 100103
$ws.Cells.Item(49, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(49, 9).Value = 100103003
$ws.Cells.Item(49, 10).Value = "Damasco"
$ws.Cells.Item(49, 11).Value = "Modesto"
$ws.Cells.Item(49, 12).Value = "Especial"
$ws.Cells.Item(49, 13).Value = 100
$ws.Cells.Item(49, 14).Value = 22000
$ws.Cells.Item(49, 15).Value = 22500
$ws.Cells.Item(49, 16).Value = 22250
$ws.Cells.Item(49, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(49, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(49, 19).Value = 1236
$ws.Cells.Item(49, 20).Value = 18

# --- New row 50: Damasco, Modesto, Primera, Región de O'Higgins, 18kg box ---
$ws.Cells.Item(50, 1).Value = 8
$ws.Cells.Item(50, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(50, 3).Value = "Coquimbo"
$ws.Cells.Item(50, 4).Value = 44568
$ws.Cells.Item(50, 5).Value = 4
$ws.Cells.Item(50, 6).Value = "Fruta"
$ws.Cells.Item(50, 7).Value = 100103
$ws.Cells.Item(50, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(50, 9).Value = 100103003
$ws.Cells.Item(50, 10).Value = "Damasco"
$ws.Cells.Item(50, 11).Value = "Modesto"
$ws.Cells.Item(50, 12).Value = "Primera"
$ws.Cells.Item(50, 13).Value = 200
$ws.Cells.Item(50, 14).Value = 20000
$ws.Cells.Item(50, 15).Value = 20500
$ws.Cells.Item(50, 16).Value = 20250
$ws.Cells.Item(50, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(50, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(50, 19).Value = 1125
$ws.Cells.Item(50, 20).Value = 18

# Make sure the date cells keep the workbook's existing datetime number format
$ws.Range("D49:D50").NumberFormat = "YYYY-MM-DD HH:MM:SS"
